$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for rows 2-8 per re-pulled data / mean calculation
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -2
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 0
